# Insert a new weekly price-report row at row 58 (pushing existing rows
# 58-93 down to 59-94), matching the "Fruta / hortaliza, semanal" update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 58:93 down one row, creating a fresh (blank) row 58.
$ws.Rows.Item(58).Insert()

# Populate the new row 58 with this week's record. Columns A,B,C,E-L,T
# mirror the constant values shared by every row in this homogeneous
# "Vega Monumental Concepción - Mango" block.
$ws.Cells.Item(58, 1).Value2 = 11
$ws.Cells.Item(58, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(58, 3).Value = "Bíobío"
$ws.Cells.Item(58, 4).Value2 = 44596
$ws.Cells.Item(58, 5).Value2 = 8
$ws.Cells.Item(58, 6).Value = "Fruta"
$ws.Cells.Item(58, 7).Value2 = 100108
$ws.Cells.Item(58, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(58, 9).Value2 = 100108002
$ws.Cells.Item(58, 10).Value = "Mango"
$ws.Cells.Item(58, 11).Value = "Sin especificar"
$ws.Cells.Item(58, 12).Value = "Primera"
$ws.Cells.Item(58, 13).Value2 = 170
$ws.Cells.Item(58, 14).Value2 = 6500
$ws.Cells.Item(58, 15).Value2 = 7000
$ws.Cells.Item(58, 16).Value2 = 6735
$ws.Cells.Item(58, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(58, 18).Value = "Ecuador"
$ws.Cells.Item(58, 19).Value2 = 1684
$ws.Cells.Item(58, 20).Value2 = 4
